# Apply the "Average/Worst of SW,SC ratios" summary additions to Sheet1.
#
# Adds:
#   J12  =AVERAGE(J2:J11)                         (bold, matches header font)
#   A14  "Average of SW(S*)/SW(OPT)"   B14 =AVERAGE(N2:N11)
#   A15  "Average of SC(S*)/SC(OPT)"   B15 =AVERAGE(Z2:Z11)
#   A16  "Worst of SW(S*)/SW(OPT)"     B16 =MIN(N2:N11)
#   A17  "Worst of SC(S*)/SC(OPT)"     B17 =MAX(Z2:Z11)
# and selects J12, and sets the page setup (paper size / orientation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: quick average of the sample-ratio column (J), bold ---
$j12 = $ws.Range("J12")
$j12.Font.Bold = $true
$j12.Formula = "=AVERAGE(J2:J11)"

# --- Labels in column A ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"

# --- Format B14 once (bold, size 12, vertically centered), then reuse the
#     same look for B15:B17 via copy/paste-format so we don't accumulate
#     extra half-built style combinations in the stylesheet. ---
$b14 = $ws.Range("B14")
$b14.VerticalAlignment = -4108
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.Formula = "=AVERAGE(N2:N11)"

$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)

$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Match the selection left behind in the saved workbook.
$ws.Range("J12").Select()

# Page setup as configured before export/print.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
